$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 270.04544
$ws.Range("I33").Value = 232.05
$ws.Range("J33").Value = 650
$ws.Range("K33").Value = 232.05
$ws.Range("L33").Value = 650
$ws.Range("M33").Value = -3.050000000000011
$ws.Range("N33").Value = -1108
$ws.Range("H40").Value = 975
$ws.Range("I40").Value = 930
$ws.Range("K40").Value = 930
$ws.Range("M40").Value = -755
$ws.Range("H51").Value = 4363.636
$ws.Range("I51").Value = 3333.3333
$ws.Range("J51").Value = 4750
$ws.Range("K51").Value = 3333.3333
$ws.Range("L51").Value = 4750
$ws.Range("M51").Value = -2849.3333
$ws.Range("N51").Value = -5718
$ws.Range("H107").Value = 882.6
$ws.Range("I107").Value = 785.43475
$ws.Range("K107").Value = 785.43475
$ws.Range("M107").Value = 1134.56525
$ws.Range("H132").Value = 202941.02
$ws.Range("I132").Value = 235810.2
$ws.Range("J132").Value = 1030.2858
$ws.Range("K132").Value = 707430.6000000001
$ws.Range("L132").Value = 3090.8574
$ws.Range("M132").Value = -704900.6000000001
$ws.Range("N132").Value = -8150.857400000001
$ws.Range("H138").Value = 2353.1636
$ws.Range("I138").Value = 1923.4546
$ws.Range("J138").Value = 2997.7273
$ws.Range("K138").Value = 5770.3638
$ws.Range("L138").Value = 8993.1819
$ws.Range("M138").Value = -630.3638000000001
$ws.Range("N138").Value = -19273.1819
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1549553.6
$ws.Range("I2").Value = 1513
$ws.Range("J2").Value = 2452577.5
$ws.Range("K2").Value = 1513
$ws.Range("L2").Value = 2452577.5
$ws.Range("M2").Value = -1400
$ws.Range("N2").Value = -2452803.5
$ws.Range("H61").Value = 12823653
$ws.Range("I61").Value = 14495954
$ws.Range("K61").Value = 14495954
$ws.Range("M61").Value = -14495742
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null
$ws.Range("H116").Value = 1549553.6
$ws.Range("I116").Value = 1513
$ws.Range("J116").Value = 2452577.5
$ws.Range("K116").Value = 1513
$ws.Range("L116").Value = 2452577.5
$ws.Range("M116").Value = 781
$ws.Range("N116").Value = -2457165.5
$ws.Range("H136").Value = 12823653
$ws.Range("I136").Value = 14495954
$ws.Range("K136").Value = 43487862
$ws.Range("M136").Value = -43485312
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1549553.6
$ws.Range("I3").Value = 1513
$ws.Range("J3").Value = 2452577.5
$ws.Range("K3").Value = 1513
$ws.Range("L3").Value = 2452577.5
$ws.Range("M3").Value = -1399
$ws.Range("N3").Value = -2452805.5
$ws.Range("H62").Value = 50000
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51372
$ws.Range("H65").Value = 50000
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156864
$ws.Range("H134").Value = 49047.44
$ws.Range("I134").Value = 55539.637
$ws.Range("J134").Value = 1438
$ws.Range("K134").Value = 166618.911
$ws.Range("L134").Value = 4314
$ws.Range("M134").Value = -164083.911
$ws.Range("N134").Value = -9384
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3862.4
$ws.Range("I16").Value = 3944.4285
$ws.Range("J16").Value = 3671
$ws.Range("K16").Value = 3944.4285
$ws.Range("L16").Value = 3671
$ws.Range("M16").Value = -3657.4285
$ws.Range("N16").Value = -4245
$ws.Range("H113").Value = 3862.4
$ws.Range("I113").Value = 3944.4285
$ws.Range("J113").Value = 3671
$ws.Range("K113").Value = 3944.4285
$ws.Range("L113").Value = 3671
$ws.Range("M113").Value = -1774.4285
$ws.Range("N113").Value = -8011
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 70707200
$ws.Range("J96").Value = 70707200
$ws.Range("L96").Value = 212121600
$ws.Range("N96").Value = -212125718
$ws.Range("H110").Value = 2839.8
$ws.Range("I110").Value = 2839.8
$ws.Range("K110").Value = 8519.400000000001
$ws.Range("M110").Value = -4429.400000000001
$ws.Range("H131").Value = 3090.6978
$ws.Range("J131").Value = 2398.7878
$ws.Range("L131").Value = 7196.3634
$ws.Range("N131").Value = -17276.3634
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1688
$ws.Range("I113").Value = 1342.6154
$ws.Range("J113").Value = 1968.625
$ws.Range("K113").Value = 1342.6154
$ws.Range("L113").Value = 1968.625
$ws.Range("M113").Value = 827.3846000000001
$ws.Range("N113").Value = -6308.625
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 549
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 498
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 498
$ws.Range("M22").Value = -305
$ws.Range("N22").Value = -1088
$ws.Range("H27").Value = 549
$ws.Range("I27").Value = 600
$ws.Range("J27").Value = 498
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 498
$ws.Range("M27").Value = -493
$ws.Range("N27").Value = -712
$ws.Range("H61").Value = 2147.7778
$ws.Range("I61").Value = 2275.7144
$ws.Range("J61").Value = 1700
$ws.Range("K61").Value = 2275.7144
$ws.Range("L61").Value = 1700
$ws.Range("M61").Value = -2073.7144
$ws.Range("N61").Value = -2104
$ws.Range("H113").Value = 2147.7778
$ws.Range("I113").Value = 2275.7144
$ws.Range("J113").Value = 1700
$ws.Range("K113").Value = 2275.7144
$ws.Range("L113").Value = 1700
$ws.Range("M113").Value = -105.7143999999998
$ws.Range("N113").Value = -6040
$ws.Range("H132").Value = 38598.1
$ws.Range("I132").Value = 44005.383
$ws.Range("J132").Value = 3450.75
$ws.Range("K132").Value = 132016.149
$ws.Range("L132").Value = 10352.25
$ws.Range("M132").Value = -129486.149
$ws.Range("N132").Value = -15412.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 348.03333
$ws.Range("I113").Value = 315.57895
$ws.Range("J113").Value = 404.0909
$ws.Range("K113").Value = 946.73685
$ws.Range("L113").Value = 1212.2727
$ws.Range("M113").Value = 1223.26315
$ws.Range("N113").Value = -5552.2727
